# 2017-2-23 MsmqTest project commit
# Update task plan sheet: mark "MVC后台管理登录实例" as completed (add end date +
# status + project name), flesh out the ".net WCF服务实例" row with full
# completed details, and start the ".net 消息队列实例" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7 first: copy the "开始" (in-progress / orange) look from E4 while
# it still carries that formatting, before we flip E4 over to "已完成".
$ws.Range("E4").Copy() | Out-Null
$ws.Range("E7").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Row 4: "MVC后台管理登录实例" is now finished.
# End date (copy the date-formatted look from an existing date cell).
$ws.Range("C4").Copy() | Out-Null
$ws.Range("D4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("D4").Value = 42786

# Completion status (copy the "已完成" green look from E2).
$ws.Range("E2").Copy() | Out-Null
$ws.Range("E4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("E4").Value = "已完成"

# Project name.
$ws.Range("F4").Value = "MVCAuthorizeTest"

# --- Row 6: ".net WCF服务实例" is complete too.
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("C6").Value = 42787

$ws.Range("D2").Copy() | Out-Null
$ws.Range("D6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("D6").Value = 42788

$ws.Range("E2").Copy() | Out-Null
$ws.Range("E6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E6").Value = "已完成"

$ws.Range("F6").Value = "WcfServiceTest"

# --- Row 7: ".net 消息队列实例" has just started.
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("C7").Value = 42789

$ws.Range("E7").Value = "开始"

# --- Match the saved selection in the workbook.
$ws.Range("D8").Select() | Out-Null
